$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 2001100524
$ws.Range("P6").Value = 2001100524
$ws.Range("P7").Value = 2001100520
$ws.Range("P8").Value = 2001100520
